$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "展览" (sheet1): update a few counts and append a new row
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F3").Value = 5474
$ws1.Range("F6").Value = 80
$ws1.Range("F8").Value = 50
$ws1.Range("F9").Value = 522

# New row 10 - copy formatting/style from row 9 first, then set values
$ws1.Range("A9:I9").Copy($ws1.Range("A10:I10"))

$ws1.Range("A10").Value = 9

# B10 holds a date-like string; force text formatting so it is not
# auto-converted into a date serial number, then restore the original
# (default) cell style by pasting formats from B9 back over it.
$ws1.Range("B10").NumberFormat = "@"
$ws1.Range("B10").Value = "2024-11-02"
$ws1.Range("B9").Copy()
$ws1.Range("B10").PasteSpecial(-4122)

$ws1.Range("C10").Value = "南宁·梦中礼Lolita茶会"
$ws1.Range("D10").Value = "吉兴西路盛天汇一、三、四层 云庭汇·安吉宴会厅"
$ws1.Range("E10").Value = "2024.11.02 13:00-11.02 17:00"
$ws1.Range("F10").Value = 3
$ws1.Range("G10").Value = 118
$ws1.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=92826"
$ws1.Range("I10").Value = "//i2.hdslb.com/bfs/openplatform/202409/09AXaAJA1726816540668.jpeg"

# ---------------------------------------------------------------
# Sheet "全部类型" (sheet4): update a few counts and append a new row
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F3").Value = 5474
$ws4.Range("F7").Value = 80
$ws4.Range("F10").Value = 50
$ws4.Range("F11").Value = 522

# New row 12 - copy formatting/style from row 11 first, then set values
$ws4.Range("A11:I11").Copy($ws4.Range("A12:I12"))

$ws4.Range("A12").Value = 11

# B12 holds a date-like string; same treatment as B10 above.
$ws4.Range("B12").NumberFormat = "@"
$ws4.Range("B12").Value = "2024-11-02"
$ws4.Range("B11").Copy()
$ws4.Range("B12").PasteSpecial(-4122)

$ws4.Range("C12").Value = "南宁·梦中礼Lolita茶会"
$ws4.Range("D12").Value = "吉兴西路盛天汇一、三、四层 云庭汇·安吉宴会厅"
$ws4.Range("E12").Value = "2024.11.02 13:00-11.02 17:00"
$ws4.Range("F12").Value = 3
$ws4.Range("G12").Value = 118
$ws4.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=92826"
$ws4.Range("I12").Value = "//i2.hdslb.com/bfs/openplatform/202409/09AXaAJA1726816540668.jpeg"
